$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "tout le monde"

$ws.Range("C42").Value = "Samira"
$ws.Range("C43").Value = "Samira"
$ws.Range("C44").Value = "Samira"
$ws.Range("C45").Value = "Samira"
$ws.Range("C46").Value = "Samira"
$ws.Range("C47").Value = "Samira"

$ws.Range("A7").Select()
